$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Contenu du stage" pie-chart source data (rows 16-23: D=language,
# E=number of students, G=percentage label as text, e.g. "7.84 %").
# Only the counts/percentages actually changed; the language labels in
# column D stay identical.

$ws.Range("E16").Value = 4
$ws.Range("E17").Value = 41
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3

# The percentage cells (column G) are stored as literal text (e.g. "0 %"),
# not as real percentage numbers. Typing a string like "7.84 %" straight
# into a General-formatted cell makes Excel auto-convert it into a number
# with a percent format, which would change the cell's style - unlike the
# original workbook where these are plain text shared strings with no
# cell style at all. Using a scratch formula cell ="text" and pasting its
# *value only* (xlPasteValues) keeps the destination a clean literal-text
# cell, matching the source workbook's layout.
$scratch = $ws.Range("ZZ1")

function Set-TextValue($targetCell, [string]$text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("G16") "7.84 %"
Set-TextValue $ws.Range("G17") "80.39 %"
Set-TextValue $ws.Range("G19") "5.88 %"
Set-TextValue $ws.Range("G20") "5.88 %"

$scratch.Clear()
